$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Data values for I2:J75
$iValues = @(7,7,6,6,8,8,6,7,8,7,8,8,7,7,7,9,6,6,6,4,6,9,7,4,1,7,1,6,7,5,7,8,6,6,1,6,7,7,5,6,8,5,6,6,6,6,8,6,6,3,7,7,8,7,7,9,6,6,4,6,8,3,8,4,7,6,6,5,6,9,6,5,8,6)
$jValues = @(7,7,7,6,8,8,7,8,8,7,8,8,7,8,7,9,7,6,8,5,7,9,7,6,2,7,1,6,7,5,7,8,6,6,2,6,7,8,6,6,8,6,7,6,7,7,8,6,7,4,8,7,8,7,8,9,7,8,5,7,8,5,8,6,7,6,7,5,7,9,6,6,8,6)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
